# Applies the diff: trims the intro/boilerplate paragraphs, rewrites the
# patch-status / RMF-compliance / next-steps / risk-assessment sections with
# the new narrative text, and removes the paragraphs that no longer exist.
#
# All paragraphs in this document are single visual "lines": the text runs
# end with a manual line break (<w:br/>, represented in Range.Text as the
# vertical-tab char, code 11) immediately before the paragraph mark
# (code 13) -- except the very last paragraph in the document, which has
# only the paragraph mark. Set-ParaText below replaces a paragraph's
# visible text while preserving that trailing break/no-break shape.

function Set-ParaText($doc, $para, $newText, $keepBreak) {
    $full = $para.Range
    # Keep the final paragraph mark (CR) in place; replace everything
    # before it (ordinarily just the pre-existing VT break char, if any).
    $r = $doc.Range($full.Start, $full.End - 1)
    if ($keepBreak) {
        $r.Text = $newText + [char]11
    } else {
        $r.Text = $newText
    }
}

$d = $word.ActiveDocument

# --- "*** Risk Assessment ***" section (originally paragraphs 32-38) ---
# Last paragraph in the doc: replace its text with the new closing "Note:" line.
Set-ParaText $d $d.Paragraphs.Item($d.Paragraphs.Count) "Note: Based on available CVE information, there are several vulnerabilities identified in various products that affect this system. However, specific details about these vulnerabilities have not been provided." $false

# Remove the now-obsolete risk/mitigation paragraphs (old 34-38), highest index first.
$d.Paragraphs.Item(38).Range.Delete()
$d.Paragraphs.Item(37).Range.Delete()
$d.Paragraphs.Item(36).Range.Delete()
$d.Paragraphs.Item(35).Range.Delete()
$d.Paragraphs.Item(34).Range.Delete()

# Insert the new risk paragraph right after the "*** Risk Assessment ***" heading.
$riskHeading = $d.Paragraphs.Item(32)
$riskHeading.Range.InsertParagraphAfter()
Set-ParaText $d $d.Paragraphs.Item(33) "There is a potential risk associated with not applying the pending patch. The impact level of this risk could be significant, including unauthorized access or exploitation of the system. To mitigate this risk, it is essential to apply the patch and maintain regular vulnerability checks to ensure the system remains secure." $true

# Fix the heading spacing: "Risk Assessment***" -> "Risk Assessment ***"
Set-ParaText $d $d.Paragraphs.Item(32) "*** Risk Assessment ***" $true

# --- "*** Recommended next steps ***" section (originally paragraphs 26-30) ---
Set-ParaText $d $d.Paragraphs.Item(30) "- Update documentation to reflect any changes made to the system as a result of the patch installation." $true
$d.Paragraphs.Item(29).Range.Delete()
Set-ParaText $d $d.Paragraphs.Item(28) "- Schedule the deployment of the patch." $true
Set-ParaText $d $d.Paragraphs.Item(27) "- Provide a review and assessment of the available patch." $true
Set-ParaText $d $d.Paragraphs.Item(26) "The recommended next steps are:" $true

# --- "*** Compliance with RMF Controls ***" section (originally paragraphs 17-23) ---
Set-ParaText $d $d.Paragraphs.Item(23) "In order to ensure compliance with Risk Management Framework (RMF) controls, it is essential to remediate any identified vulnerabilities. This involves identifying the vulnerabilities in place, reporting them, and taking corrective action to address them. Additionally, configuration management should be implemented to track changes and updates made to the system. Regular vulnerability checks should also be performed to ensure the system remains secure." $true
$d.Paragraphs.Item(22).Range.Delete()
$d.Paragraphs.Item(21).Range.Delete()
$d.Paragraphs.Item(20).Range.Delete()
$d.Paragraphs.Item(19).Range.Delete()
$d.Paragraphs.Item(18).Range.Delete()
$d.Paragraphs.Item(17).Range.Delete()

# --- "*** Patch Status Summary ***" section (originally paragraphs 13-14) ---
Set-ParaText $d $d.Paragraphs.Item(14) "The system is currently running with pending patches. There is a single patch available, which is code/stable 1.99.0-1743632463 amd64 [upgradable from: 1.98.2-1741788907]. This patch is relevant to security as it addresses vulnerabilities that could allow unauthorized access or exploitation of the system." $true
$d.Paragraphs.Item(13).Range.Delete()

# --- "*** System Overview ***" section (originally paragraphs 3-4) ---
$d.Paragraphs.Item(4).Range.Delete()
$d.Paragraphs.Item(3).Range.Delete()
